$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "WallDrop" row (row 31) under the existing table, matching the
# formatting used by the nearby rows (numeric-only display format).
$ws.Cells.Item(31, 1).Value = "WallDrop"
$ws.Cells.Item(31, 2).Value = 1
$ws.Cells.Item(31, 3).Value = 12

$newRow = $ws.Range("A31:C31")
$newRow.NumberFormat = "0"

# Update the view: scroll so row 16 is at the top and select cell B33,
# matching the author's final cursor position.
$ws.Range("B33").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
